# Clear column C "DietaryRestrictions2" values for rows 2-31 (the "None"
# placeholder column), keeping their existing style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = $null
}

# Scroll the view down and update the active selection to C32:C62.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C32:C62").Select()
